$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week dates) ---
$ws.Cells.Item(8, 1).Value = "Volume 31   Number  16"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# --- Reference cells for style-preserving type changes (never modified themselves) ---
$refZero  = $ws.Cells.Item(14, 3)   # style 14, text "0"  (C14)
$refNA    = $ws.Cells.Item(14, 5)   # style 14, text "***.*" (E14)
$refNum15 = $ws.Cells.Item(15, 9)   # style 15, numeric (I15, unused as data target)
$refNum16 = $ws.Cells.Item(15, 11)  # style 16, numeric (K15, unused as data target)

# --- Cells changing type: fix style+type via Copy from reference, then set exact value ---
$refZero.Copy($ws.Cells.Item(15, 4))
# D15 -> text "0" (no further value set needed)

$refNA.Copy($ws.Cells.Item(15, 5))
# E15 -> text "***.*" (no further value set needed)

$refZero.Copy($ws.Cells.Item(16, 3))
# C16 -> text "0" (no further value set needed)

$refZero.Copy($ws.Cells.Item(18, 3))
# C18 -> text "0" (no further value set needed)

$refZero.Copy($ws.Cells.Item(22, 3))
# C22 -> text "0" (no further value set needed)

$refZero.Copy($ws.Cells.Item(27, 4))
# D27 -> text "0" (no further value set needed)

$refNA.Copy($ws.Cells.Item(27, 5))
# E27 -> text "***.*" (no further value set needed)

$refNum15.Copy($ws.Cells.Item(28, 3))
$ws.Cells.Item(28, 3).Value = 1

$refNum15.Copy($ws.Cells.Item(33, 4))
$ws.Cells.Item(33, 4).Value = 1

$refNum16.Copy($ws.Cells.Item(33, 5))
$ws.Cells.Item(33, 5).Value = -100

$refZero.Copy($ws.Cells.Item(33, 6))
# F33 -> text "0" (no further value set needed)

# --- Plain value updates (type unchanged) ---
$ws.Cells.Item(15, 14).Value = 75
$ws.Cells.Item(16, 4).Value = 2
$ws.Cells.Item(16, 5).Value = -100
$ws.Cells.Item(16, 6).Value = 8
$ws.Cells.Item(16, 7).Value = 7
$ws.Cells.Item(16, 8).Value = 14.285714285714
$ws.Cells.Item(16, 10).Value = 25
$ws.Cells.Item(16, 11).Value = 8
$ws.Cells.Item(16, 13).Value = -49.056603773584
$ws.Cells.Item(16, 14).Value = -88.260869565217
$ws.Cells.Item(17, 3).Value = 7
$ws.Cells.Item(17, 4).Value = 3
$ws.Cells.Item(17, 5).Value = 133.333333333333
$ws.Cells.Item(17, 6).Value = 24
$ws.Cells.Item(17, 7).Value = 8
$ws.Cells.Item(17, 8).Value = 200
$ws.Cells.Item(17, 9).Value = 67
$ws.Cells.Item(17, 10).Value = 55
$ws.Cells.Item(17, 11).Value = 21.818181818181
$ws.Cells.Item(17, 12).Value = 13.559322033898
$ws.Cells.Item(17, 13).Value = 71.794871794871
$ws.Cells.Item(17, 14).Value = -11.842105263157
$ws.Cells.Item(18, 4).Value = 5
$ws.Cells.Item(18, 5).Value = -100
$ws.Cells.Item(18, 6).Value = 13
$ws.Cells.Item(18, 7).Value = 11
$ws.Cells.Item(18, 8).Value = 18.181818181818
$ws.Cells.Item(18, 10).Value = 41
$ws.Cells.Item(18, 11).Value = -7.317073170731
$ws.Cells.Item(18, 12).Value = -28.301886792452
$ws.Cells.Item(18, 13).Value = -67.241379310344
$ws.Cells.Item(18, 14).Value = -92.910447761194
$ws.Cells.Item(19, 3).Value = 10
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 5).Value = 11.111111111111
$ws.Cells.Item(19, 6).Value = 45
$ws.Cells.Item(19, 8).Value = 36.363636363636
$ws.Cells.Item(19, 9).Value = 174
$ws.Cells.Item(19, 10).Value = 163
$ws.Cells.Item(19, 11).Value = 6.748466257668
$ws.Cells.Item(19, 12).Value = 10.828025477707
$ws.Cells.Item(19, 13).Value = 50
$ws.Cells.Item(19, 14).Value = -13
$ws.Cells.Item(20, 3).Value = 6
$ws.Cells.Item(20, 4).Value = 4
$ws.Cells.Item(20, 5).Value = 50
$ws.Cells.Item(20, 6).Value = 16
$ws.Cells.Item(20, 7).Value = 11
$ws.Cells.Item(20, 8).Value = 45.454545454545
$ws.Cells.Item(20, 9).Value = 55
$ws.Cells.Item(20, 10).Value = 44
$ws.Cells.Item(20, 11).Value = 25
$ws.Cells.Item(20, 12).Value = 111.538461538462
$ws.Cells.Item(20, 13).Value = 22.222222222222
$ws.Cells.Item(20, 14).Value = -90.848585690515
$ws.Cells.Item(21, 3).Value = 23
$ws.Cells.Item(21, 4).Value = 23
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 106
$ws.Cells.Item(21, 7).Value = 73
$ws.Cells.Item(21, 8).Value = 45.205479452054
$ws.Cells.Item(21, 9).Value = 370
$ws.Cells.Item(21, 10).Value = 335
$ws.Cells.Item(21, 11).Value = 10.447761194029
$ws.Cells.Item(21, 12).Value = 12.121212121212
$ws.Cells.Item(21, 13).Value = -0.269541778975
$ws.Cells.Item(21, 14).Value = -77.602905569007
$ws.Cells.Item(24, 3).Value = 26
$ws.Cells.Item(24, 4).Value = 24
$ws.Cells.Item(24, 5).Value = 8.333333333333
$ws.Cells.Item(24, 6).Value = 93
$ws.Cells.Item(24, 7).Value = 79
$ws.Cells.Item(24, 8).Value = 17.721518987341
$ws.Cells.Item(24, 9).Value = 310
$ws.Cells.Item(24, 10).Value = 331
$ws.Cells.Item(24, 11).Value = -6.344410876132
$ws.Cells.Item(24, 12).Value = -10.404624277456
$ws.Cells.Item(24, 13).Value = 23.505976095617
$ws.Cells.Item(25, 4).Value = 9
$ws.Cells.Item(25, 5).Value = -55.555555555555
$ws.Cells.Item(25, 6).Value = 17
$ws.Cells.Item(25, 7).Value = 26
$ws.Cells.Item(25, 8).Value = -34.615384615384
$ws.Cells.Item(25, 9).Value = 70
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = -30
$ws.Cells.Item(25, 12).Value = -36.363636363636
$ws.Cells.Item(26, 3).Value = 12
$ws.Cells.Item(26, 4).Value = 9
$ws.Cells.Item(26, 5).Value = 33.333333333333
$ws.Cells.Item(26, 6).Value = 33
$ws.Cells.Item(26, 7).Value = 26
$ws.Cells.Item(26, 8).Value = 26.923076923076
$ws.Cells.Item(26, 9).Value = 132
$ws.Cells.Item(26, 10).Value = 89
$ws.Cells.Item(26, 11).Value = 48.314606741573
$ws.Cells.Item(26, 12).Value = 37.5
$ws.Cells.Item(26, 13).Value = 6.451612903225
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(28, 5).Value = -66.666666666666
$ws.Cells.Item(28, 6).Value = 2
$ws.Cells.Item(28, 7).Value = 6
$ws.Cells.Item(28, 8).Value = -66.666666666666
$ws.Cells.Item(28, 9).Value = 24
$ws.Cells.Item(28, 10).Value = 14
$ws.Cells.Item(28, 11).Value = 71.428571428571
$ws.Cells.Item(28, 12).Value = 33.333333333333
$ws.Cells.Item(31, 12).Value = 16.666666666666
$ws.Cells.Item(33, 7).Value = 2
$ws.Cells.Item(33, 8).Value = -100
$ws.Cells.Item(33, 10).Value = 2
$ws.Cells.Item(33, 11).Value = -50
